# informe_diario_20250408.xlsx - "cambio por tema de cifras"
# Updates several percentage figures (stored as text in the sheet) and a
# handful of plain numeric totals on the "Resumen" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Percentage values (stored as literal text, e.g. "0.48%") ---------
# The workbook keeps these as plain text cells (no numeric formatting),
# so we force a Text number format before writing the new value to stop
# Excel from reinterpreting "0.48%" as the number 0.0048.
$percentCells = @(
    @("B10", "0.48%"),
    @("C10", "0.33%"),
    @("D10", "0.81%"),
    @("D11", "55.76%"),
    @("B12", "22.42%"),
    @("C12", "21.81%"),
    @("D12", "44.24%"),
    @("B15", "26.73%"),
    @("D15", "76.96%"),
    @("B16", "9.68%"),
    @("D16", "23.04%")
)

foreach ($pair in $percentCells) {
    $cell = $ws.Range($pair[0])
    $cell.NumberFormat = "@"
    $cell.Value = $pair[1]
}

# --- Plain numeric totals ----------------------------------------------
$ws.Range("D2").Value = 16996

$ws.Range("B13").Value = 58
$ws.Range("D13").Value = 167

$ws.Range("B14").Value = 21
$ws.Range("D14").Value = 50
